$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the highlight fill style from these cells (they previously used the
# "theme 9" fill style, which is being removed from the workbook entirely).
$clearFillCells = @("F3", "G3", "H3", "I3", "A4", "F5", "G5", "H5", "I5", "A6")
foreach ($addr in $clearFillCells) {
    $ws.Range($addr).Style = "Normal"
}

# New "X" entries in the extended status-word columns for rows 4 and 6.
$ws.Range("K4:O4").Value = "X"
$ws.Range("K6:O6").Value = "X"

# Update the active selection.
$ws.Range("F3:I6").Select()
$excel.ActiveCell = $ws.Range("F3")
